$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.240107999999999
$ws.Range("H2").Value = 18.720324
$ws.Range("I2").Value = 0.01732230523539376
$ws.Range("J2").Value = 0.01732230523539376
$ws.Range("M2").Value = 0.305927
$ws.Range("N2").Value = 0.917781
$ws.Range("O2").Value = 0.01078151192043695
$ws.Range("P2").Value = 0.01078151192043695
$ws.Range("Q2").Value = 1.909017520116
$ws.Range("R2").Value = 17.181157681044
$ws.Range("S2").Value = 0.0001867606403848452
$ws.Range("T2").Value = 0.0001867606403848452
$ws.Range("G3").Value = 6.240107999999999
$ws.Range("H3").Value = 18.720324
$ws.Range("I3").Value = 0.01732230523539376
$ws.Range("J3").Value = 0.01732230523539376
$ws.Range("O3").Value = 0.09577486496546363
$ws.Range("P3").Value = 0.09577486496546361
$ws.Range("Q3").Value = 16.958279743608
$ws.Range("R3").Value = 152.624517692472
$ws.Range("S3").Value = 0.001659041444810381
$ws.Range("T3").Value = 0.001659041444810381
$ws.Range("G4").Value = 6.240107999999999
$ws.Range("H4").Value = 18.720324
$ws.Range("I4").Value = 0.01732230523539376
$ws.Range("J4").Value = 0.01732230523539376
$ws.Range("O4").Value = 0.8934436231140994
$ws.Range("P4").Value = 0.8934436231140994
$ws.Range("Q4").Value = 158.196692852296
$ws.Range("R4").Value = 1423.770235670664
$ws.Range("S4").Value = 0.01547650315019853
$ws.Range("T4").Value = 0.01547650315019854
$ws.Range("I5").Value = 0.9592798330716089
$ws.Range("J5").Value = 0.9592798330716091
$ws.Range("M5").Value = 0.305927
$ws.Range("N5").Value = 0.917781
$ws.Range("O5").Value = 0.01078151192043695
$ws.Range("P5").Value = 0.01078151192043695
$ws.Range("Q5").Value = 105.718146813733
$ws.Range("R5").Value = 951.4633213235969
$ws.Range("S5").Value = 0.01034248695529632
$ws.Range("T5").Value = 0.01034248695529632
$ws.Range("I6").Value = 0.9592798330716089
$ws.Range("J6").Value = 0.9592798330716091
$ws.Range("O6").Value = 0.09577486496546363
$ws.Range("P6").Value = 0.09577486496546361
$ws.Range("S6").Value = 0.09187489647652583
$ws.Range("T6").Value = 0.09187489647652584
$ws.Range("I7").Value = 0.9592798330716089
$ws.Range("J7").Value = 0.9592798330716091
$ws.Range("O7").Value = 0.8934436231140994
$ws.Range("P7").Value = 0.8934436231140994
$ws.Range("S7").Value = 0.8570624496397867
$ws.Range("T7").Value = 0.8570624496397869
$ws.Range("G8").Value = 8.428738666666666
$ws.Range("I8").Value = 0.02339786169299727
$ws.Range("J8").Value = 0.02339786169299728
$ws.Range("M8").Value = 0.305927
$ws.Range("N8").Value = 0.917781
$ws.Range("O8").Value = 0.01078151192043695
$ws.Range("P8").Value = 0.01078151192043695
$ws.Range("Q8").Value = 2.578578734077333
$ws.Range("R8").Value = 23.207208606696
$ws.Range("S8").Value = 0.0002522643247557851
$ws.Range("T8").Value = 0.0002522643247557852
$ws.Range("G9").Value = 8.428738666666666
$ws.Range("I9").Value = 0.02339786169299727
$ws.Range("J9").Value = 0.02339786169299728
$ws.Range("O9").Value = 0.09577486496546363
$ws.Range("P9").Value = 0.09577486496546361
$ws.Range("Q9").Value = 22.90615934773867
$ws.Range("R9").Value = 206.155434129648
$ws.Range("S9").Value = 0.002240927044127408
$ws.Range("T9").Value = 0.002240927044127408
$ws.Range("G10").Value = 8.428738666666666
$ws.Range("I10").Value = 0.02339786169299727
$ws.Range("J10").Value = 0.02339786169299728
$ws.Range("O10").Value = 0.8934436231140994
$ws.Range("P10").Value = 0.8934436231140994
$ws.Range("R10").Value = 1923.137746629776
$ws.Range("S10").Value = 0.02090467032411408
$ws.Range("T10").Value = 0.02090467032411408
